# Add a new keyword variant ("<program>변경하려") to the keyword list (column D)
# for the eight "APPLY_CHANGE" FAQ rows (rows 26-33), one per program in
# column C. This mirrors the upstream edit which appended an extra
# keyword to each of these shared-string entries (the workbook's save
# process also reshuffled the shared-strings table, but that is an
# internal/cosmetic artifact of re-saving - the actual user-visible
# content change is just this appended keyword).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 26..33

foreach ($r in $rows) {
    $program = $ws.Cells.Item($r, 3).Text
    $current = $ws.Cells.Item($r, 4).Text
    $ws.Cells.Item($r, 4).Value = $current + ", " + $program + "변경하려"
}

# Match the author's final view state: scrolled down so row 16 is at the
# top, with D34 as the active/selected cell.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D34").Select()
